# "Add files via upload" — re-upload of the Internship Checklist workbook
# with a few checklist rows filled in / relabeled.
#
# Changes applied (per the target diff):
#   - D27: "Cover letter" -> "Cover"
#   - E10: (blank) -> "Ask mentor to review"
#   - E11: (blank) -> "Ask mentor to review"
#   - E29: (blank) -> "Mentor visits company"
#   - Selection moves from E31 to E32
#
# NOTE: order matters so the new shared-string table entries land in the
# same order as the target file (Cover, Ask mentor to review, Mentor visits
# company).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D27").Value = "Cover"
$ws.Range("E10").Value = "Ask mentor to review"
$ws.Range("E11").Value = "Ask mentor to review"
$ws.Range("E29").Value = "Mentor visits company"

$ws.Range("E32").Select()
